# "Generate Report for Archive"
#
# The localization status moved from "Ready for handoff" to
# "In Translation" everywhere it appears (Overview!E2:F2, zh-cn!C2,
# de-de!C2). Re-generating the report also re-sized the (now shorter)
# status column on each sheet.

$wb = $excel.ActiveWorkbook

# --- Update the status text wherever it appears -----------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Shrink the status columns to match the new (shorter) text --------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
